# Commit: 1st commit on 3rd July 2021
#
# - Rename the only sheet "USDBTC" -> "PAIR2"
# - Duplicate it to a new sheet "VNDUSD" (keeps header row + styling)
# - Update the data rows on both sheets with the new trade info

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet ------------------------------------------------
$pair2 = $wb.Worksheets.Item(1)
$pair2.Name = "PAIR2"

# --- Swap the "1HR CHART" / "15MIN CHART" header columns -------------------
$e1 = $pair2.Range("E1").Value2
$f1 = $pair2.Range("F1").Value2
$pair2.Range("E1").Value = $f1
$pair2.Range("F1").Value = $e1

# --- Update PAIR2 row 2 data -----------------------------------------------
$pair2.Range("A2").Value = 7
$pair2.Range("B2").Value = 7
$pair2.Range("C2").Value = "14:58:01.136710"
$pair2.Range("D2").Value = "Sell"
$pair2.Range("E2").Value = "link 1"
$pair2.Range("F2").Value = "link 2"
$pair2.Range("G2").Value = 3
$pair2.Range("H2").Value = "this is my comment"
$pair2.Range("I2").Value = 810
$pair2.Range("J2").Value = 3

# --- Create VNDUSD as a copy of PAIR2 (keeps headers/column widths/styles) -
$pair2.Copy([System.Reflection.Missing]::Value, $pair2)
$vndusd = $wb.Worksheets.Item(2)
$vndusd.Name = "VNDUSD"

# --- Update VNDUSD row 2 data -----------------------------------------------
$vndusd.Range("A2").Value = 7
$vndusd.Range("B2").Value = 7
$vndusd.Range("C2").Value = "14:57:30.000793"
$vndusd.Range("D2").Value = "Sell"
$vndusd.Range("E2").Value = "link 1"
$vndusd.Range("F2").Value = "link 2"
$vndusd.Range("G2").Value = 3
$vndusd.Range("H2").Value = "this is my comment"
$vndusd.Range("I2").Value = 610
$vndusd.Range("J2").Value = 3

# --- Add VNDUSD row 3 (new trade) -------------------------------------------
# Give row 3 the same cell styling as row 2 (border + centered) before filling it in
$vndusd.Range("A2:J2").Copy()
$vndusd.Range("A3:J3").PasteSpecial(-4122)  # xlPasteFormats

$vndusd.Range("A3").Value = 11
$vndusd.Range("B3").Value = 7
$vndusd.Range("C3").Value = "14:57:03.901108"
$vndusd.Range("D3").Value = "Sell"
$vndusd.Range("E3").Value = "link 1"
$vndusd.Range("F3").Value = "link 2"
$vndusd.Range("G3").Value = 3
$vndusd.Range("H3").Value = "this is my comment"
$vndusd.Range("I3").Value = 410
$vndusd.Range("J3").Value = 3

# Make sure PAIR2 is the active/selected sheet, matching the original workbook
$pair2.Activate()
